# Update the "Responsible?" column of the indicator-sources table on slide 2:
# add "/Valeria" to "Ruben" for the Corruption Perception Index (CPI),
# Government Effectiveness Indicator (GEI), Politics right Index (PRI) and
# Civil Liberties index rows.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(1)
$tbl = $shape.Table

# Column 7 = "Responsible?" ; rows 4-7 = CPI, GEI, PRI, Civil Liberties index
for ($r = 4; $r -le 7; $r++) {
    $cellRange = $tbl.Cell($r, 7).Shape.TextFrame.TextRange
    if ($cellRange.Text -eq "Ruben") {
        $cellRange.Text = "Ruben/Valeria"
    }
}
